$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: create blank template rows 777-798 (copy blank format currently on row 768) ---
$ws.Range("A768:E768").Copy()
for ($r = 777; $r -le 798; $r++) {
    $ws.Range("A" + $r + ":E" + $r).PasteSpecial(-4122)
    $ws.Rows.Item($r).RowHeight = 15
}

# --- Step 2: create data-row formatting for new rows 769-776 (copy from existing data row 767) ---
$ws.Range("A767:E767").Copy()
for ($r = 769; $r -le 776; $r++) {
    $ws.Range("A" + $r + ":E" + $r).PasteSpecial(-4122)
    $ws.Rows.Item($r).RowHeight = 15
}

# --- Step 3: fix formatting of row 768 (was blank template, needs to become a data row) ---
$ws.Range("A767:E767").Copy()
$ws.Range("A768:E768").PasteSpecial(-4122)
$ws.Rows.Item(768).RowHeight = 15

# --- Step 4: populate the 9 new event rows (768-776) ---

# Row 768
$ws.Range("A768").Value = 45920
$ws.Range("B768").Value = "BINARY ROOFTOP RAVE"
$ws.Range("C768").Value = "Innenhafen"
$ws.Range("D768").Value = "Duisburg"
$txt = "https://www.instagram.com/reel/DOY3-mNDUew/?igsh=MWVmZ3Fja2lqMnN6ZQ=="
$ws.Range("E768").Value = $txt
$ws.Hyperlinks.Add($ws.Range("E768"), $txt, "", "", $txt)
$len = $txt.Length
$c1 = $ws.Range("E768").Characters(1, $len - 1)
$c1.Font().Underline = 2
$c1.Font().Color = 65280
$c1.Font().Name = "Calibri"
$c1.Font().Size = 11
$c2 = $ws.Range("E768").Characters($len, 1)
$c2.Font().Underline = 2
$c2.Font().Color = 65280
$c2.Font().Name = "Calibri"
$c2.Font().Size = 11
$ws.Range("B767").Copy()
$ws.Range("E768").PasteSpecial(-4122)

# Row 769
$ws.Range("A769").Value = 45989
$ws.Range("B769").Value = "ZWILLING & DJ ACHIM FEUERVOGEL"
$ws.Range("C769").Value = "Artheater"
$ws.Range("D769").Value = "Köln"
$txt = "https://www.instagram.com/reel/DOlV84IDJZp/?igsh=YWJ0NThhbmI1bnBw"
$ws.Range("E769").Value = $txt
$ws.Hyperlinks.Add($ws.Range("E769"), $txt, "", "", $txt)
$len = $txt.Length
$c1 = $ws.Range("E769").Characters(1, $len - 1)
$c1.Font().Underline = 2
$c1.Font().Color = 65280
$c1.Font().Name = "Calibri"
$c1.Font().Size = 11
$c2 = $ws.Range("E769").Characters($len, 1)
$c2.Font().Underline = 2
$c2.Font().Color = 65280
$c2.Font().Name = "Calibri"
$c2.Font().Size = 11
$ws.Range("B767").Copy()
$ws.Range("E769").PasteSpecial(-4122)

# Row 770
$ws.Range("A770").Value = 45920
$ws.Range("B770").Value = "SONIC ESCAPE"
$ws.Range("C770").Value = "PM93"
$ws.Range("D770").Value = "Essen"
$txt = "https://www.instagram.com/reel/DOodUdZjIuM/?igsh=a3F6bG95a3NjaHdz"
$ws.Range("E770").Value = $txt
$ws.Hyperlinks.Add($ws.Range("E770"), $txt, "", "", $txt)
$len = $txt.Length
$c1 = $ws.Range("E770").Characters(1, $len - 1)
$c1.Font().Underline = 2
$c1.Font().Color = 65280
$c1.Font().Name = "Calibri"
$c1.Font().Size = 11
$c2 = $ws.Range("E770").Characters($len, 1)
$c2.Font().Underline = 2
$c2.Font().Color = 65280
$c2.Font().Name = "Calibri"
$c2.Font().Size = 11
$ws.Range("B767").Copy()
$ws.Range("E770").PasteSpecial(-4122)

# Row 771
$ws.Range("A771").Value = 45919
$ws.Range("B771").Value = "BASSMANIA"
$ws.Range("C771").Value = "Favela"
$ws.Range("D771").Value = "Münster"
$txt = "https://www.instagram.com/reel/DOosRHACFyr/?igsh=MWR2cDNjZmtldjl4cA=="
$ws.Range("E771").Value = $txt
$ws.Hyperlinks.Add($ws.Range("E771"), $txt, "", "", $txt)
$len = $txt.Length
$c1 = $ws.Range("E771").Characters(1, $len - 1)
$c1.Font().Underline = 2
$c1.Font().Color = 65280
$c1.Font().Name = "Calibri"
$c1.Font().Size = 11
$c2 = $ws.Range("E771").Characters($len, 1)
$c2.Font().Underline = 2
$c2.Font().Color = 65280
$c2.Font().Name = "Calibri"
$c2.Font().Size = 11
$ws.Range("B767").Copy()
$ws.Range("E771").PasteSpecial(-4122)

# Row 772
$ws.Range("A772").Value = 45927
$ws.Range("B772").Value = "CLUB NIGHT"
$ws.Range("C772").Value = "SNRS"
$ws.Range("D772").Value = "Dortmund"
$txt = "https://www.instagram.com/reel/DObA3FcjMUN/?igsh=ZmF0ZTBmYjk4aGZx"
$ws.Range("E772").Value = $txt
$ws.Hyperlinks.Add($ws.Range("E772"), $txt, "", "", $txt)
$len = $txt.Length
$c1 = $ws.Range("E772").Characters(1, $len - 1)
$c1.Font().Underline = 2
$c1.Font().Color = 65280
$c1.Font().Name = "Calibri"
$c1.Font().Size = 11
$c2 = $ws.Range("E772").Characters($len, 1)
$c2.Font().Underline = 2
$c2.Font().Color = 65280
$c2.Font().Name = "Calibri"
$c2.Font().Size = 11
$ws.Range("B767").Copy()
$ws.Range("E772").PasteSpecial(-4122)

# Row 773
$ws.Range("A773").Value = 45969
$ws.Range("B773").Value = "MATTERMIND"
$ws.Range("C773").Value = "Essigfabrik & Elektroküche"
$ws.Range("D773").Value = "Köln"
$txt = "https://www.instagram.com/reel/DOn5TaEDLUJ/?igsh=aXY0OXB4ZXpqMGw1"
$ws.Range("E773").Value = $txt
$ws.Hyperlinks.Add($ws.Range("E773"), $txt, "", "", $txt)
$len = $txt.Length
$c1 = $ws.Range("E773").Characters(1, $len - 1)
$c1.Font().Underline = 2
$c1.Font().Color = 65280
$c1.Font().Name = "Calibri"
$c1.Font().Size = 11
$c2 = $ws.Range("E773").Characters($len, 1)
$c2.Font().Underline = 2
$c2.Font().Color = 65280
$c2.Font().Name = "Calibri"
$c2.Font().Size = 11
$ws.Range("B767").Copy()
$ws.Range("E773").PasteSpecial(-4122)

# Row 774
$ws.Range("A774").Value = 45932
$ws.Range("B774").Value = "WYLDHEARTS"
$ws.Range("C774").Value = "Helios37"
$ws.Range("D774").Value = "Köln"
$txt = "https://www.instagram.com/reel/DOoYF_0j7dl/?igsh=cG5scDR3cnJmbWZo"
$ws.Range("E774").Value = $txt
$ws.Hyperlinks.Add($ws.Range("E774"), $txt, "", "", $txt)
$len = $txt.Length
$c1 = $ws.Range("E774").Characters(1, $len - 1)
$c1.Font().Underline = 2
$c1.Font().Color = 65280
$c1.Font().Name = "Calibri"
$c1.Font().Size = 11
$c2 = $ws.Range("E774").Characters($len, 1)
$c2.Font().Underline = 2
$c2.Font().Color = 65280
$c2.Font().Name = "Calibri"
$c2.Font().Size = 11
$ws.Range("B767").Copy()
$ws.Range("E774").PasteSpecial(-4122)

# Row 775
$ws.Range("A775").Value = 45927
$ws.Range("B775").Value = "TECHNO EVOLUTION EVENTS"
$ws.Range("C775").Value = "Am Hawerkamp 31"
$ws.Range("D775").Value = "Münster"
$txt = "https://www.instagram.com/reel/DN3kQW10LT5/?igsh=MXRrMDZiY2dqdG5rdw=="
$ws.Range("E775").Value = $txt
$ws.Hyperlinks.Add($ws.Range("E775"), $txt, "", "", $txt)
$len = $txt.Length
$c1 = $ws.Range("E775").Characters(1, $len - 1)
$c1.Font().Underline = 2
$c1.Font().Color = 65280
$c1.Font().Name = "Calibri"
$c1.Font().Size = 11
$c2 = $ws.Range("E775").Characters($len, 1)
$c2.Font().Underline = 2
$c2.Font().Color = 65280
$c2.Font().Name = "Calibri"
$c2.Font().Size = 11
$ws.Range("B767").Copy()
$ws.Range("E775").PasteSpecial(-4122)

# Row 776
$ws.Range("A776").Value = 45919
$ws.Range("B776").Value = "SOLI-PARTY"
$ws.Range("C776").Value = "Fade‘in"
$ws.Range("D776").Value = "Münster"
$txt = "https://www.instagram.com/reel/DNTRzeUCai_/?igsh=MTk1cTZ2cTgzMmh5eA=="
$ws.Range("E776").Value = $txt
$ws.Hyperlinks.Add($ws.Range("E776"), $txt, "", "", $txt)
$len = $txt.Length
$c1 = $ws.Range("E776").Characters(1, $len - 1)
$c1.Font().Underline = 2
$c1.Font().Color = 65280
$c1.Font().Name = "Calibri"
$c1.Font().Size = 11
$c2 = $ws.Range("E776").Characters($len, 1)
$c2.Font().Underline = 2
$c2.Font().Color = 65280
$c2.Font().Name = "Calibri"
$c2.Font().Size = 11
$ws.Range("B767").Copy()
$ws.Range("E776").PasteSpecial(-4122)

Write-Host "All rows populated."
